$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New names for rows 2-11 (column B)
$names = @(
    "Juliana da Mata",
    "Sr. Gustavo Henrique da Cunha",
    "Lucas Freitas",
    "Dra. Helena da Conceição",
    "Sr. Heitor Lima",
    "Sr. Marcos Vinicius da Cruz",
    "Rebeca Campos",
    "Elisa Gomes",
    "Emanuella Cardoso",
    "Elisa da Cunha"
)

# New Faltas values (column C)
$faltas = @(1, 8, 5, 3, 20, 1, 1, 4, 5, 11)

# New Nota values (column D)
$notas = @(80, 50, 75, 50, 40, 95, 85, 50, 55, 75)

for ($i = 0; $i -lt 10; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $names[$i]
    $ws.Cells.Item($row, 3).Value = $faltas[$i]
    $ws.Cells.Item($row, 4).Value = $notas[$i]
}
